$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.244.14"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").Value = "1.900.64"
$ws.Range("E3").Value = "  -0.25%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").Value = "'326.29"
$ws.Range("E5").Value = "  -0.55%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  -0.26%  "
$ws.Range("D7").Value = "'0.4642"
$ws.Range("D8").Value = "'0.3915"
$ws.Range("E8").Value = "  -1.04%  "
$ws.Range("D9").Value = "'0.07890"
$ws.Range("E9").Value = "  -0.92%  "
$ws.Range("E10").Value = "  -0.99%  "
$ws.Range("D11").Value = "'21.80"
$ws.Range("E11").Value = "  -2.15%  "
$ws.Range("D12").Value = "1.961.27"
$ws.Range("E12").Value = "  +1.71%  "
$ws.Range("D13").Value = "'7.074"
$ws.Range("E13").Value = "  -0.82%  "
$ws.Range("D14").Value = "'5.743"
$ws.Range("E14").Value = "  -0.43%  "
$ws.Range("D15").Value = "'0.07003"
$ws.Range("E15").Value = "  +0.45%  "
$ws.Range("D16").Value = "'88.24"
$ws.Range("E16").Value = "  -0.55%  "
$ws.Range("D17").Value = "'1.002"
$ws.Range("E17").Value = "  -0.28%  "
$ws.Range("D18").Value = "'0.000009984"
$ws.Range("E18").Value = "  -1.28%  "
$ws.Range("D19").Value = "'17.11"
$ws.Range("E19").Value = "  -0.27%  "
$ws.Range("D21").Value = "29.260.14"
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("D22").Value = "'5.296"
$ws.Range("E22").Value = "  -1.09%  "
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").Value = "2.149.01"
$ws.Range("E24").Value = "  -0.34%  "
$ws.Range("D25").Value = "'2.100"
$ws.Range("E25").Value = "  +2.46%  "
$ws.Range("D26").Value = "'156.23"
$ws.Range("E26").Value = "  -0.48%  "
$ws.Range("D27").Value = "'19.45"
$ws.Range("E27").Value = "  -0.28%  "
$ws.Range("D28").Value = "'5.979"
$ws.Range("E28").Value = "  +1.33%  "
$ws.Range("D29").Value = "'118.71"
$ws.Range("E29").Value = "  -0.46%  "
$ws.Range("D30").Value = "'1.886"
$ws.Range("E30").Value = "  -5.57%  "
$ws.Range("D31").Value = "'0.09331"
$ws.Range("E31").Value = "  -0.87%  "
$ws.Range("D32").Value = "'0.9009"
$ws.Range("E32").Value = "  -2.38%  "
$ws.Range("D33").Value = "'5.266"
$ws.Range("E33").Value = "  -1.55%  "
$ws.Range("D34").Value = "'1.326"
$ws.Range("E34").Value = "  -1.55%  "
$ws.Range("E35").Value = "  -1.57%  "
$ws.Range("D36").Value = "'1.190"
$ws.Range("E36").Value = "  +1.73%  "
$ws.Range("D37").Value = "'0.05777"
$ws.Range("E37").Value = "  -0.88%  "
$ws.Range("E38").Value = "  -0.83%  "
$ws.Range("E39").Value = "  -0.26%  "
$ws.Range("D40").Value = "'7.712"
$ws.Range("E40").Value = "  -3.59%  "
$ws.Range("E41").Value = "  -0.84%  "
$ws.Range("D42").Value = "'0.1788"
$ws.Range("E42").Value = "  -1.25%  "
$ws.Range("D43").Value = "'9.708"
$ws.Range("E43").Value = "  -2.72%  "
$ws.Range("D44").Value = "'11.98"
$ws.Range("E44").Value = "  -0.58%  "
$ws.Range("D45").Value = "'0.5358"
$ws.Range("E45").Value = "  -1.30%  "
$ws.Range("D46").Value = "'2.173"
$ws.Range("E46").Value = "  -2.03%  "
$ws.Range("D47").Value = "'0.07020"
$ws.Range("E47").Value = "  -1.14%  "
$ws.Range("D48").Value = "'1.853"
$ws.Range("E48").Value = "  -1.26%  "
$ws.Range("D49").Value = "'2.574"
$ws.Range("E49").Value = "  -0.94%  "
$ws.Range("D50").Value = "'113.17"
$ws.Range("E50").Value = "  +1.03%  "
$ws.Range("D51").Value = "'1.060"
$ws.Range("E51").Value = "  -0.34%  "
